# msz - mandatory fields checks part 1
# Add a new "102_AutomobileInsurance_002_VehicleData_001_MandatoryFields" record
# row to the flow table, mirroring the pattern of the existing
# "Open Automobile Insurance" (row 3) entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "102_AutomobileInsurance_002_VehicleData_001_MandatoryFields"
$ws.Range("B7").Value = "var102_AutomobileInsurance_002_VehicleData_001_MandatoryFields"
$ws.Range("C7").Value = "Open Automobile Insurance"
$ws.Range("D7").Value = "102_AutomobileInsurance_002_VehicleData_001_MandatoryFields"

# New text is longer than what used to fit in columns A, B and D, so
# Excel re-autosizes them to fit the new longest entries.
$ws.Columns("A").ColumnWidth = 55.666666666666664
$ws.Columns("B").ColumnWidth = 58.333333333333336
$ws.Columns("D").ColumnWidth = 55.666666666666664

# Move / reflect the current selection like the saved workbook shows.
$ws.Range("D8").Select()

# Reposition/resize the workbook window to match the saved view state.
$win = $excel.ActiveWindow
$win.Left = 3852
$win.Top = 3444
$win.Width = 29964
$win.Height = 9180
